$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Enterprise Grid Security" title textbox (TextBox 4) ---
$titleBox = $s.Shapes.Item("TextBox 4")

# Nudge the textbox down (y offset 152400 -> 228600 EMU, i.e. 12pt -> 18pt).
$titleBox.Top = 18

# The title used to be two runs ("Enterprise Grid " + "Security") with a
# leftover endParaRPr on the paragraph. Clear the whole range first (an
# explicit Delete(), not just Text = "") so retyping the text rebuilds the
# paragraph as a single run with no trailing endParaRPr, then set the
# combined text.
$titleBox.TextFrame.TextRange.Delete()
$titleBox.TextFrame.TextRange.Text = "Enterprise Grid Security"

# --- "http://www.cagrid.org" textbox (TextBox 9) ---
$cagridBox = $s.Shapes.Item("TextBox 9")
$cagridBox.TextFrame.TextRange.Delete()
$cagridBox.TextFrame.TextRange.Text = "http://www.cagrid.org"

# --- "GAARDS Security Infrastructure" textbox (TextBox 10) ---
$gaardsBox = $s.Shapes.Item("TextBox 10")
$gaardsBox.TextFrame.TextRange.Delete()
$gaardsBox.TextFrame.TextRange.Text = "GAARDS Security Infrastructure"
